$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "gruenenthal" / "theramex" rows (53 & 54) ---
# Row 53 now refers to theramex, row 54 now refers to gruenenthal.
$ws.Range("B53").Value = "theramex"
$ws.Range("B54").Value = "gruenenthal"

# Row 53 (theramex) new figures
$ws.Range("A53").Value = 53
$ws.Range("C53").Value = 26
$ws.Range("D53").Value = 61065
$ws.Range("E53").Value = 2349

# Row 54 (gruenenthal) new figures
$ws.Range("A54").Value = 23
$ws.Range("C54").Value = 153
$ws.Range("D54").Value = 329993
$ws.Range("E54").Value = 2157

# --- Reruns / numeric fixes across the sheet ---
$ws.Range("D4").Value = 16715103
$ws.Range("E4").Value = 40181

$ws.Range("D10").Value = 4537484
$ws.Range("E10").Value = 26228

$ws.Range("D16").Value = 3119493.91

$ws.Range("D17").Value = 11915637.61

$ws.Range("D22").Value = 386943.46

$ws.Range("D24").Value = 2467592

$ws.Range("D26").Value = 3759699.3

$ws.Range("D31").Value = 1317573.67

$ws.Range("D32").Value = 2517371.2

$ws.Range("D34").Value = 1785331.25
$ws.Range("E34").Value = 5282

$ws.Range("D40").Value = 216392

$ws.Range("D50").Value = 173693.91

$ws.Range("D56").Value = 118542.02
$ws.Range("E56").Value = 1428
